$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "About" sheet (sheet 1): refresh the "last updated" date and make it the
#    active/selected tab (it previously lived on "Current and Planned
#    Capacity").
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value2 = 45387

# ---------------------------------------------------------------------------
# 2) "BAU Emissions" sheet (sheet 5): bulk rename the "NoSettings" model-run
#    label to "test" across every row label in column A, update the cached
#    run values for row 94 (natural gas / iron & steel), and move the
#    viewport/selection to match the author's latest scroll position.
# ---------------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$wsBau.Cells.Replace(" : NoSettings", " : test") | Out-Null

$wsBau.Range("M94").Value2 = 1001080
$wsBau.Range("N94").Value2 = 2002150
$wsBau.Range("O94").Value2 = 3003230
$wsBau.Range("P94").Value2 = 4004300
$wsBau.Range("Q94").Value2 = 5005380
$wsBau.Range("R94").Value2 = 5005380
$wsBau.Range("S94").Value2 = 5005380
$wsBau.Range("T94").Value2 = 5005380
$wsBau.Range("U94").Value2 = 5005380
$wsBau.Range("V94").Value2 = 5005380
$wsBau.Range("W94").Value2 = 5005380
$wsBau.Range("X94").Value2 = 5005380
$wsBau.Range("Y94").Value2 = 5005380
$wsBau.Range("Z94").Value2 = 5005380
$wsBau.Range("AA94").Value2 = 5005380
$wsBau.Range("AB94").Value2 = 5005380
$wsBau.Range("AC94").Value2 = 5005380
$wsBau.Range("AD94").Value2 = 5005380
$wsBau.Range("AE94").Value2 = 5005380

# Move the BAU Emissions viewport/selection before we leave it, then make
# "About" the active sheet/tab as the workbook should open there.
$wsBau.Activate()
$wsBau.Range("A30:AE280").Select()

$wsAbout.Activate()
$wsAbout.Range("E29").Select()
